$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# P2 / P3 hold tracking numbers that are stored as text (shared strings) in
# the workbook, even though they look like numbers. Force a text entry
# (equivalent to typing a leading apostrophe in Excel) so the values are
# not silently converted to numeric cells.
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "320018207767"

$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "320018207778"
